# EvaluacionPorCompetenciasTemplate.xlsx edit
# Commit: "Mostrar promedio ponderado, y nombre de meritos y competencias en el reporte"
#
# Changes:
#  1. Remove the "Tipo de Reclutamiento:" and "Origen de la Evaluación:" form fields,
#     shifting "Nombre Evaluador:" / "Fecha de Evaluación:" up.
#  2. Clear the extra custom-competency rows (37-43: Prueba Técnica Requisito Adicional 1/2,
#     Solución de problemas, Negociación, Planificación, Trabajar bajo presión,
#     Tomar decisiones en el momento) leaving only Capacitación Formal/Complementaria/Experiencia.
#  3. Update both charts so they stop referencing the removed rows/series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- 1. Evaluator header block (rows 4-7) ---------------------------------
# Before:
#   C4 Nombre:                     H4 Cédula:
#   C5 Nombre del puesto:          H5 Tipo de Reclutamiento:
#   C6 Origen de la Evaluación:    H6 Nombre Evaluador:
#   C7 (input)                     H7 Fecha de Evaluación:
# After:
#   C4 Nombre:                     H4 Cédula:
#   C5 Nombre del puesto:          H5 Nombre Evaluador:
#   C6 (blank)                     H6 Fecha de Evaluación:
#   C7 (input)                     (H7/I7 removed)
$ws.Range("H5").Value = "Nombre Evaluador:"
$ws.Range("H6").Value = "Fecha de Evaluación:"
$ws.Range("C6").Value = ""
$ws.Range("H7:I7").UnMerge()
$ws.Range("H7").Value = ""
$ws.Range("I7").Value = ""

# --- 2. Clear the extra competency rows (37-43) ----------------------------
$ws.Range("E37:F43").ClearContents()
$ws.Range("G37:I43").ClearContents()

# Row 44 (blank spacer row) loses its E:F merge in the edited template.
$ws.Range("E44:F44").UnMerge()

# --- 3. Update chart1 (bar chart, "Calificación de Competencias") --------
$chart1 = $ws.ChartObjects(1).Chart
# Drop the trailing series that referenced row 44 (always-blank spacer row).
$chart1.SeriesCollection(11).Delete()

# --- 4. Update chart2 (radar chart, "Candidato Ideal" / "Calificación") --
$chart2 = $ws.ChartObjects(2).Chart
$chart2.SeriesCollection(1).Formula = "=SERIES(Hoja1!`$H`$33,Hoja1!`$E`$34:`$E`$43,Hoja1!`$H`$34:`$H`$43,1)"
$chart2.SeriesCollection(2).Formula = "=SERIES(Hoja1!`$I`$33,Hoja1!`$E`$34:`$E`$43,Hoja1!`$I`$34:`$I`$43,2)"

# --- 5. Restore the view/selection to match the edited area ---------------
$ws.Range("E43:I43").Select()
